$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 387
$ws.Range("I33").Value = 156.85715
$ws.Range("K33").Value = 156.85715
$ws.Range("M33").Value = 72.14285000000001

$ws.Range("H51").Value = 6677.7856
$ws.Range("I51").Value = 4462.7334
$ws.Range("J51").Value = 9233.615
$ws.Range("K51").Value = 4462.7334
$ws.Range("L51").Value = 9233.615
$ws.Range("M51").Value = -3978.7334
$ws.Range("N51").Value = -10201.615

$ws.Range("H62").Value = 2471.2083
$ws.Range("I62").Value = 2366.8096
$ws.Range("J62").Value = 3202
$ws.Range("K62").Value = 2366.8096
$ws.Range("L62").Value = 3202
$ws.Range("M62").Value = -1742.8096
$ws.Range("N62").Value = -4450

$ws.Range("H65").Value = 2471.2083
$ws.Range("I65").Value = 2366.8096
$ws.Range("J65").Value = 3202
$ws.Range("K65").Value = 11834.048
$ws.Range("L65").Value = 16010
$ws.Range("M65").Value = -8714.048000000001
$ws.Range("N65").Value = -22250

$ws.Range("H95").Value = 29998.5
$ws.Range("J95").Value = 29998.5
$ws.Range("L95").Value = 29998.5
$ws.Range("N95").Value = -35490.5

$ws.Range("H125").Value = 6283.615
$ws.Range("I125").Value = 1733.3334
$ws.Range("J125").Value = 7648.7
$ws.Range("K125").Value = 15600.0006
$ws.Range("L125").Value = 68838.3
$ws.Range("M125").Value = -13140.0006
$ws.Range("N125").Value = -73758.3

$ws.Range("H136").Value = 118593.336
$ws.Range("J136").Value = 118593.336
$ws.Range("L136").Value = 118593.336
$ws.Range("N136").Value = -128793.336

$ws.Range("H137").Value = 7307524.5
$ws.Range("I137").Value = 323884.84
$ws.Range("K137").Value = 971654.52
$ws.Range("M137").Value = -969104.52

$ws.Range("H138").Value = 3146.7476
$ws.Range("I138").Value = 1379.3077
$ws.Range("J138").Value = 3776.2466
$ws.Range("K138").Value = 4137.9231
$ws.Range("L138").Value = 11328.7398
$ws.Range("M138").Value = 1002.0769
$ws.Range("N138").Value = -21608.7398

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23419.6
$ws.Range("I32").Value = 24754.549
$ws.Range("J32").Value = 18821.445
$ws.Range("K32").Value = 24754.549
$ws.Range("L32").Value = 18821.445
$ws.Range("M32").Value = -24467.549
$ws.Range("N32").Value = -19395.445

$ws.Range("H45").Value = 2797.5454
$ws.Range("I45").Value = 1295.5
$ws.Range("K45").Value = 1295.5
$ws.Range("M45").Value = -918.5

$ws.Range("H97").Value = 613.06665
$ws.Range("I97").Value = 649.6923
$ws.Range("J97").Value = 375
$ws.Range("K97").Value = 649.6923
$ws.Range("L97").Value = 375
$ws.Range("M97").Value = -153.6923
$ws.Range("N97").Value = -1367

$ws.Range("H110").Value = 1363568.2
$ws.Range("I110").Value = 1459608.9
$ws.Range("K110").Value = 1459608.9
$ws.Range("M110").Value = -1457563.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3044.742
$ws.Range("I20").Value = 2626.05
$ws.Range("J20").Value = 3806
$ws.Range("K20").Value = 2626.05
$ws.Range("L20").Value = 3806
$ws.Range("M20").Value = -2379.05
$ws.Range("N20").Value = -4300

$ws.Range("H86").Value = 7199.143
$ws.Range("I86").Value = 5873.75
$ws.Range("K86").Value = 5873.75
$ws.Range("M86").Value = -4750.75

$ws.Range("H89").Value = 7199.143
$ws.Range("I89").Value = 5873.75
$ws.Range("K89").Value = 29368.75
$ws.Range("M89").Value = -23752.75

$ws.Range("H105").Value = 2400
$ws.Range("I105").Value = 2400
$ws.Range("K105").Value = 2400
$ws.Range("M105").Value = -653

$ws.Range("H134").Value = 1510.3928
$ws.Range("I134").Value = 1464
$ws.Range("J134").Value = 1788.75
$ws.Range("K134").Value = 4392
$ws.Range("L134").Value = 5366.25
$ws.Range("M134").Value = -1857
$ws.Range("N134").Value = -10436.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 324.6154
$ws.Range("I7").Value = 269
$ws.Range("J7").Value = 359.375
$ws.Range("K7").Value = 269
$ws.Range("L7").Value = 359.375
$ws.Range("M7").Value = -156
$ws.Range("N7").Value = -585.375

$ws.Range("H59").Value = 49999.5
$ws.Range("J59").Value = 49999.5
$ws.Range("L59").Value = 49999.5
$ws.Range("N59").Value = -52289.5

$ws.Range("H60").Value = 45611.11
$ws.Range("J60").Value = 67870.164
$ws.Range("L60").Value = 67870.164
$ws.Range("N60").Value = -68892.164

$ws.Range("H62").Value = 65442.2
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 80552.75
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 80552.75
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -81800.75

$ws.Range("H65").Value = 65442.2
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 80552.75
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 402763.75
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -409003.75

$ws.Range("H132").Value = 7585916.5
$ws.Range("I132").Value = 8341843
$ws.Range("K132").Value = 25025529
$ws.Range("M132").Value = -25022999

$ws.Range("H134").Value = 1266.5902
$ws.Range("I134").Value = 1267.7255
$ws.Range("K134").Value = 3803.1765
$ws.Range("M134").Value = -1268.1765

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12945
$ws.Range("I3").Value = 890
$ws.Range("J3").Value = 25000
$ws.Range("K3").Value = 2670
$ws.Range("L3").Value = 75000
$ws.Range("M3").Value = -2558
$ws.Range("N3").Value = -75224

$ws.Range("H98").Value = 1149
$ws.Range("J98").Value = 1235.5333
$ws.Range("L98").Value = 3706.5999
$ws.Range("N98").Value = -6702.5999

$ws.Range("H122").Value = 1299.9
$ws.Range("J122").Value = 1517.1428
$ws.Range("L122").Value = 13654.2852
$ws.Range("N122").Value = -18554.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 3859
$ws.Range("I33").Value = 2140
$ws.Range("J33").Value = 5291.5
$ws.Range("K33").Value = 2140
$ws.Range("L33").Value = 5291.5
$ws.Range("M33").Value = -1888
$ws.Range("N33").Value = -5795.5

$ws.Range("H80").Value = 843775.3
$ws.Range("I80").Value = 1388883.6
$ws.Range("J80").Value = 26112.875
$ws.Range("K80").Value = 1388883.6
$ws.Range("L80").Value = 26112.875
$ws.Range("M80").Value = -1387885.6
$ws.Range("N80").Value = -28108.875

$ws.Range("H83").Value = 843775.3
$ws.Range("I83").Value = 1388883.6
$ws.Range("J83").Value = 26112.875
$ws.Range("K83").Value = 6944418
$ws.Range("L83").Value = 130564.375
$ws.Range("M83").Value = -6939426
$ws.Range("N83").Value = -140548.375

$ws.Range("H122").Value = 382291.8
$ws.Range("I122").Value = 580559.6
$ws.Range("K122").Value = 1741678.8
$ws.Range("M122").Value = -1739228.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6607.769
$ws.Range("I7").Value = 4650.1665
$ws.Range("K7").Value = 4650.1665
$ws.Range("M7").Value = -4538.1665

$ws.Range("H22").Value = 1154.4546
$ws.Range("I22").Value = 999.8
$ws.Range("J22").Value = 1283.3334
$ws.Range("K22").Value = 999.8
$ws.Range("L22").Value = 1283.3334
$ws.Range("M22").Value = -704.8
$ws.Range("N22").Value = -1873.3334

$ws.Range("H27").Value = 1154.4546
$ws.Range("I27").Value = 999.8
$ws.Range("J27").Value = 1283.3334
$ws.Range("K27").Value = 999.8
$ws.Range("L27").Value = 1283.3334
$ws.Range("M27").Value = -892.8
$ws.Range("N27").Value = -1497.3334

$ws.Range("H40").Value = 3954.4119
$ws.Range("I40").Value = 3954.4119
$ws.Range("K40").Value = 3954.4119
$ws.Range("M40").Value = -3818.4119

$ws.Range("H46").Value = 5858.4546
$ws.Range("I46").Value = 3898
$ws.Range("J46").Value = 6386.269
$ws.Range("K46").Value = 3898
$ws.Range("L46").Value = 6386.269
$ws.Range("M46").Value = -3710
$ws.Range("N46").Value = -6762.269

$ws.Range("H55").Value = 52631880
$ws.Range("I55").Value = 125000170
$ws.Range("J55").Value = 396.18182
$ws.Range("K55").Value = 125000170
$ws.Range("L55").Value = 396.18182
$ws.Range("M55").Value = -124999997
$ws.Range("N55").Value = -742.18182

$ws.Range("H126").Value = 6607.769
$ws.Range("I126").Value = 4650.1665
$ws.Range("K126").Value = 13950.4995
$ws.Range("M126").Value = -11480.4995

$ws.Range("H132").Value = 3751.6584
$ws.Range("I132").Value = 3652.2354
$ws.Range("J132").Value = 4234.5713
$ws.Range("K132").Value = 10956.7062
$ws.Range("L132").Value = 12703.7139
$ws.Range("M132").Value = -8426.706200000001
$ws.Range("N132").Value = -17763.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6700.3335
$ws.Range("I96").Value = 1295.5
$ws.Range("J96").Value = 9402.75
$ws.Range("K96").Value = 1295.5
$ws.Range("L96").Value = 9402.75
$ws.Range("M96").Value = 77.5
$ws.Range("N96").Value = -12148.75

$ws.Range("H136").Value = 6375.96
$ws.Range("I136").Value = 2317.3462
$ws.Range("J136").Value = 10772.792
$ws.Range("K136").Value = 6952.0386
$ws.Range("L136").Value = 32318.376
$ws.Range("M136").Value = -4402.0386
$ws.Range("N136").Value = -37418.376
